$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3525
$ws.Range("I74").Value = 3020
$ws.Range("J74").Value = 4366.6665
$ws.Range("K74").Value = 3020
$ws.Range("L74").Value = 4366.6665
$ws.Range("M74").Value = -2084
$ws.Range("N74").Value = -6238.6665

$ws.Range("H77").Value = 3525
$ws.Range("I77").Value = 3020
$ws.Range("J77").Value = 4366.6665
$ws.Range("K77").Value = 15100
$ws.Range("L77").Value = 21833.3325
$ws.Range("M77").Value = -10420
$ws.Range("N77").Value = -31193.3325

$ws.Range("H86").Value = 4000.7144
$ws.Range("I86").Value = 3000.25
$ws.Range("J86").Value = 5334.6665
$ws.Range("K86").Value = 3000.25
$ws.Range("L86").Value = 5334.6665
$ws.Range("M86").Value = -1877.25
$ws.Range("N86").Value = -7580.6665

$ws.Range("H89").Value = 4000.7144
$ws.Range("I89").Value = 3000.25
$ws.Range("J89").Value = 5334.6665
$ws.Range("K89").Value = 15001.25
$ws.Range("L89").Value = 26673.3325
$ws.Range("M89").Value = -9385.25
$ws.Range("N89").Value = -37905.3325

$ws.Range("H107").Value = 1070.9
$ws.Range("I107").Value = 989.8889
$ws.Range("K107").Value = 989.8889
$ws.Range("M107").Value = 930.1111

$ws.Range("H116").Value = 1905.5385
$ws.Range("I116").Value = 2224.4
$ws.Range("J116").Value = 1706.25
$ws.Range("K116").Value = 2224.4
$ws.Range("L116").Value = 1706.25
$ws.Range("M116").Value = 1217.6
$ws.Range("N116").Value = -8590.25

$ws.Range("H129").Value = 7578.1333
$ws.Range("I129").Value = 17275.834
$ws.Range("J129").Value = 1113
$ws.Range("K129").Value = 51827.50199999999
$ws.Range("L129").Value = 3339
$ws.Range("M129").Value = -46827.50199999999
$ws.Range("N129").Value = -13339

$ws.Range("H132").Value = 4037034
$ws.Range("I132").Value = 4722066
$ws.Range("J132").Value = 2956.2222
$ws.Range("K132").Value = 14166198
$ws.Range("L132").Value = 8868.6666
$ws.Range("M132").Value = -14163668
$ws.Range("N132").Value = -13928.6666

$ws.Range("H138").Value = 3609.8235
$ws.Range("I138").Value = 3253.4167
$ws.Range("J138").Value = 3686.1965
$ws.Range("K138").Value = 9760.250100000001
$ws.Range("L138").Value = 11058.5895
$ws.Range("M138").Value = -4620.250100000001
$ws.Range("N138").Value = -21338.5895

$ws.Range("H139").Value = 65880
$ws.Range("J139").Value = 65880
$ws.Range("L139").Value = 65880
$ws.Range("N139").Value = -76160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33382.57
$ws.Range("I32").Value = 15155.1045
$ws.Range("J32").Value = 63552.17
$ws.Range("K32").Value = 15155.1045
$ws.Range("L32").Value = 63552.17
$ws.Range("M32").Value = -14868.1045
$ws.Range("N32").Value = -64126.17

$ws.Range("H61").Value = 1852.0968
$ws.Range("I61").Value = 1523.1052
$ws.Range("K61").Value = 1523.1052
$ws.Range("M61").Value = -1311.1052

$ws.Range("H80").Value = 25676.666
$ws.Range("J80").Value = 25676.666
$ws.Range("L80").Value = 25676.666
$ws.Range("N80").Value = -27672.666

$ws.Range("H83").Value = 25676.666
$ws.Range("J83").Value = 25676.666
$ws.Range("L83").Value = 77029.99800000001
$ws.Range("N83").Value = -87013.99800000001

$ws.Range("H132").Value = 19382.455
$ws.Range("I132").Value = 21507.31
$ws.Range("J132").Value = 3977.25
$ws.Range("K132").Value = 64521.93000000001
$ws.Range("L132").Value = 11931.75
$ws.Range("M132").Value = -61991.93000000001
$ws.Range("N132").Value = -16991.75

$ws.Range("H136").Value = 1852.0968
$ws.Range("I136").Value = 1523.1052
$ws.Range("K136").Value = 4569.3156
$ws.Range("M136").Value = -2019.3156

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H105").Value = 335418.66
$ws.Range("I105").Value = 252597.5
$ws.Range("J105").Value = 501061
$ws.Range("K105").Value = 252597.5
$ws.Range("L105").Value = 501061
$ws.Range("M105").Value = -250850.5
$ws.Range("N105").Value = -504555

$ws.Range("H134").Value = 3898.2886
$ws.Range("I134").Value = 3921.3953
$ws.Range("J134").Value = 3787.889
$ws.Range("K134").Value = 11764.1859
$ws.Range("L134").Value = 11363.667
$ws.Range("M134").Value = -9229.1859
$ws.Range("N134").Value = -16433.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2134.1482
$ws.Range("I58").Value = 1959.2941
$ws.Range("J58").Value = 2431.4
$ws.Range("K58").Value = 1959.2941
$ws.Range("L58").Value = 2431.4
$ws.Range("M58").Value = -1756.2941
$ws.Range("N58").Value = -2837.4

$ws.Range("H105").Value = 1348
$ws.Range("I105").Value = 1338.6666
$ws.Range("K105").Value = 1338.6666
$ws.Range("M105").Value = 408.3334

$ws.Range("H134").Value = 1216.2333
$ws.Range("I134").Value = 702.6316
$ws.Range("J134").Value = 2103.3635
$ws.Range("K134").Value = 2107.8948
$ws.Range("L134").Value = 6310.0905
$ws.Range("M134").Value = 427.1052
$ws.Range("N134").Value = -11380.0905

$ws.Range("H136").Value = 2134.1482
$ws.Range("I136").Value = 1959.2941
$ws.Range("J136").Value = 2431.4
$ws.Range("K136").Value = 5877.8823
$ws.Range("L136").Value = 7294.200000000001
$ws.Range("M136").Value = -3327.8823
$ws.Range("N136").Value = -12394.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 39902.832
$ws.Range("J137").Value = 8781.315000000001
$ws.Range("L137").Value = 26343.945
$ws.Range("N137").Value = -36543.945

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2742.818
$ws.Range("I132").Value = 2161.6155
$ws.Range("J132").Value = 3582.3333
$ws.Range("K132").Value = 6484.8465
$ws.Range("L132").Value = 10746.9999
$ws.Range("M132").Value = -3954.8465
$ws.Range("N132").Value = -15806.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9011.429
$ws.Range("I46").Value = 8170
$ws.Range("J46").Value = 10133.333
$ws.Range("K46").Value = 8170
$ws.Range("L46").Value = 10133.333
$ws.Range("M46").Value = -7982
$ws.Range("N46").Value = -10509.333

$ws.Range("H94").Value = 31077.5
$ws.Range("J94").Value = 31077.5
$ws.Range("L94").Value = 31077.5
$ws.Range("N94").Value = -32429.5

$ws.Range("H123").Value = 26681.818
$ws.Range("J123").Value = 26681.818
$ws.Range("L123").Value = 26681.818
$ws.Range("N123").Value = -36481.818

$ws.Range("H132").Value = 4503.92
$ws.Range("I132").Value = 5905.385
$ws.Range("J132").Value = 2985.6667
$ws.Range("K132").Value = 17716.155
$ws.Range("L132").Value = 8957.000100000001
$ws.Range("M132").Value = -15186.155
$ws.Range("N132").Value = -14017.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H62").Value = 15388234
$ws.Range("I62").Value = 38466536
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 38466536
$ws.Range("L62").Value = 2700
$ws.Range("M62").Value = -38465912
$ws.Range("N62").Value = -3948

$ws.Range("H65").Value = 15388234
$ws.Range("I65").Value = 38466536
$ws.Range("J65").Value = 2700
$ws.Range("K65").Value = 192332680
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -192329560
$ws.Range("N65").Value = -19740

$ws.Range("H122").Value = 2125.6924
$ws.Range("I122").Value = 2219.5
$ws.Range("K122").Value = 6658.5
$ws.Range("M122").Value = -4208.5

$ws.Range("H123").Value = 17238.096
$ws.Range("J123").Value = 17238.096
$ws.Range("L123").Value = 17238.096
$ws.Range("N123").Value = -27038.096

$ws.Range("H132").Value = 26372.762
$ws.Range("I132").Value = 2705.5312
$ws.Range("J132").Value = 102107.9
$ws.Range("K132").Value = 8116.5936
$ws.Range("L132").Value = 306323.7
$ws.Range("M132").Value = -5586.5936
$ws.Range("N132").Value = -311383.7

$ws.Range("H136").Value = 2804.4412
$ws.Range("I136").Value = 3522.4707
$ws.Range("J136").Value = 2086.4119
$ws.Range("K136").Value = 10567.4121
$ws.Range("L136").Value = 6259.2357
$ws.Range("M136").Value = -8017.4121
$ws.Range("N136").Value = -11359.2357
